$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = 44312
$ws.Range("J2").Value2 = 180
$ws.Range("D3").Value2 = 44497
$ws.Range("I3").Value2 = 'Primera'
$ws.Range("J3").Value2 = 250
$ws.Range("K3").Value2 = 800
$ws.Range("L3").Value2 = 800
$ws.Range("M3").Value2 = 800
$ws.Range("N3").Value2 = '$/kilo (volumen en unidades)'
$ws.Range("O3").Value2 = 'Perú'
$ws.Range("P3").Value2 = 800
$ws.Range("D4").Value2 = 44495
$ws.Range("K4").Value2 = 800
$ws.Range("L4").Value2 = 800
$ws.Range("M4").Value2 = 800
$ws.Range("N4").Value2 = '$/kilo (volumen en unidades)'
$ws.Range("O4").Value2 = 'Perú'
$ws.Range("P4").Value2 = 800
$ws.Range("D5").Value2 = 44223
$ws.Range("H5").Value2 = 'Americana O Klondike'
$ws.Range("I5").Value2 = 'Extra'
$ws.Range("J5").Value2 = 340
$ws.Range("K5").Value2 = 2500
$ws.Range("L5").Value2 = 2500
$ws.Range("M5").Value2 = 2500
$ws.Range("P5").Value2 = 2500
$ws.Range("D6").Value2 = 44223
$ws.Range("H6").Value2 = 'Americana O Klondike'
$ws.Range("I6").Value2 = 'Primera'
$ws.Range("J6").Value2 = 400
$ws.Range("K6").Value2 = 2000
$ws.Range("L6").Value2 = 2000
$ws.Range("M6").Value2 = 2000
$ws.Range("P6").Value2 = 2000
$ws.Range("D7").Value2 = 44223
$ws.Range("H7").Value2 = 'Americana O Klondike'
$ws.Range("I7").Value2 = 'Segunda'
$ws.Range("J7").Value2 = 300
$ws.Range("K7").Value2 = 1500
$ws.Range("L7").Value2 = 1500
$ws.Range("M7").Value2 = 1500
$ws.Range("P7").Value2 = 1500
$ws.Range("D8").Value2 = 44223
$ws.Range("H8").Value2 = 'Americana O Klondike'
$ws.Range("I8").Value2 = 'Tercera'
$ws.Range("J8").Value2 = 160
$ws.Range("K8").Value2 = 1000
$ws.Range("L8").Value2 = 1000
$ws.Range("M8").Value2 = 1000
$ws.Range("P8").Value2 = 1000
$ws.Range("D9").Value2 = 44167
$ws.Range("J9").Value2 = 400
$ws.Range("K9").Value2 = 5000
$ws.Range("L9").Value2 = 5000
$ws.Range("M9").Value2 = 5000
$ws.Range("P9").Value2 = 5000
$ws.Range("D10").Value2 = 44167
$ws.Range("I10").Value2 = 'Segunda'
$ws.Range("J10").Value2 = 560
$ws.Range("K10").Value2 = 3000
$ws.Range("L10").Value2 = 3000
$ws.Range("M10").Value2 = 3000
$ws.Range("N10").Value2 = '$/unidad'
$ws.Range("O10").Value2 = 'Región de O''Higgins'
$ws.Range("P10").Value2 = 3000
$ws.Range("D11").Value2 = 44167
$ws.Range("I11").Value2 = 'Tercera'
$ws.Range("J11").Value2 = 450
$ws.Range("K11").Value2 = 2000
$ws.Range("L11").Value2 = 2000
$ws.Range("M11").Value2 = 2000
$ws.Range("N11").Value2 = '$/unidad'
$ws.Range("O11").Value2 = 'Región de O''Higgins'
$ws.Range("P11").Value2 = 2000
$ws.Range("D12").Value2 = 44488
$ws.Range("J12").Value2 = 150
$ws.Range("D14").Value2 = 44194
$ws.Range("I14").Value2 = 'Extra'
$ws.Range("J14").Value2 = 120
$ws.Range("K14").Value2 = 3500
$ws.Range("L14").Value2 = 3500
$ws.Range("M14").Value2 = 3500
$ws.Range("N14").Value2 = '$/unidad'
$ws.Range("O14").Value2 = 'Región de O''Higgins'
$ws.Range("P14").Value2 = 3500
$ws.Range("D15").Value2 = 44194
$ws.Range("J15").Value2 = 200
$ws.Range("K15").Value2 = 3000
$ws.Range("L15").Value2 = 3000
$ws.Range("M15").Value2 = 3000
$ws.Range("N15").Value2 = '$/unidad'
$ws.Range("O15").Value2 = 'Región de O''Higgins'
$ws.Range("P15").Value2 = 3000
$ws.Range("D16").Value2 = 44305
$ws.Range("H16").Value2 = 'Sin especificar'
$ws.Range("I16").Value2 = 'Primera'
$ws.Range("J16").Value2 = 100
$ws.Range("O16").Value2 = 'Perú'
$ws.Range("D17").Value2 = 44477
$ws.Range("H17").Value2 = 'Sin especificar'
$ws.Range("J17").Value2 = 80
$ws.Range("K17").Value2 = 800
$ws.Range("L17").Value2 = 800
$ws.Range("M17").Value2 = 800
$ws.Range("N17").Value2 = '$/kilo (volumen en unidades)'
$ws.Range("O17").Value2 = 'Perú'
$ws.Range("P17").Value2 = 800
$ws.Range("D18").Value2 = 44217
$ws.Range("H18").Value2 = 'Sin especificar'
$ws.Range("I18").Value2 = 'Extra'
$ws.Range("J18").Value2 = 400
$ws.Range("K18").Value2 = 2500
$ws.Range("L18").Value2 = 2500
$ws.Range("M18").Value2 = 2500
$ws.Range("P18").Value2 = 2500
$ws.Range("D19").Value2 = 44217
$ws.Range("H19").Value2 = 'Sin especificar'
$ws.Range("I19").Value2 = 'Primera'
$ws.Range("J19").Value2 = 280
$ws.Range("K19").Value2 = 2000
$ws.Range("L19").Value2 = 2000
$ws.Range("M19").Value2 = 2000
$ws.Range("P19").Value2 = 2000
$ws.Range("D20").Value2 = 44491
$ws.Range("J20").Value2 = 150
$ws.Range("K20").Value2 = 800
$ws.Range("L20").Value2 = 800
$ws.Range("M20").Value2 = 800
$ws.Range("N20").Value2 = '$/kilo (volumen en unidades)'
$ws.Range("P20").Value2 = 800
